# Update hotel reviews data: populate the previously-blank
# English_Reviews_num (G2) / Local_Rank (H2) / Total_Reviews_num (I2)
# cells on the hotel_info sheet with the scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hotel_info")

# Leading apostrophe forces these numeric-looking values to be stored
# as text (matching the source data's shared-string type), then we
# reset the style so no extra number-format/quote-prefix style sticks
# to the cell.
$ws.Range("G2").Value = "'5"
$ws.Range("G2").Style = "Normal"

$ws.Range("H2").Value = "'419"
$ws.Range("H2").Style = "Normal"

$ws.Range("I2").Value = "'5"
$ws.Range("I2").Style = "Normal"
